# Label BOM items better.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resistor package label: "R-W4" -> "R-1/4W"
$ws.Range("C2").Value = "R-1/4W"
$ws.Range("C3").Value = "R-1/4W"

# Ceramic capacitor package label: "C-5mm" -> "C-P5mm"
$ws.Range("C4").Value = "C-P5mm"

# Electrolytic capacitor package label: "E2.5-6.3" -> "E-P2.5mm 6.3x11.5mm"
$ws.Range("C5").Value = "E-P2.5mm 6.3x11.5mm"

# Ceramic capacitor description: "Capacitor Ceramic THT" -> "Ceramic Capacitor THT"
$ws.Range("E4").Value = "Ceramic Capacitor THT"

# Electrolytic capacitor description: "Capacitor Polarized THT" -> "Electrolytic Capacitor THT"
$ws.Range("E5").Value = "Electrolytic Capacitor THT"

# Update the active selection to match the authored state
$ws.Range("E4:E5").Select()
